$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: "No" -> "Yes"
$ws.Range("B2").Value = "Yes"

# F2: empty -> "Na"
$ws.Range("F2").Value = "Na"

# J2: 0 -> 1
$ws.Range("J2").Value = 1

# F3: "FISHER" -> "Fisher"
$ws.Range("F3").Value = "Fisher"
